$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("馬際泰(c_personid=677027)的入仕方式是？", "舉人")
    ,@("夏錫金(c_personid=344371)的入仕方式是？", "進士")
    ,@("鍾德浦(c_personid=667675)的入仕方式是？", "恩蔭")
    ,@("悟勤善淳(c_personid=375984)的入仕方式是？", "進士")
    ,@("陳錡(c_personid=124113的籍貫是否為郵州？", "否")
    ,@("王禹功(c_personid=342115)的入仕方式是？", "進士")
    ,@("展毓(c_personid=198547的籍貫是否為騎右衛？", "否")
    ,@("潘之驤(c_personid=634068)的籍貫是？", "仁和")
    ,@("盧庭琮(c_personid=348921的籍貫是否為國州？", "否")
    ,@("程天秩(c_personid=7135)的入仕方式是？", "恩蔭")
    ,@("吳洵武(c_personid=27040)的入仕方式是？", "舉人")
    ,@("陳錡(c_personid=124113的籍貫是否為揚州府？", "是")
    ,@("謝顯光(c_personid=658370)的入仕方式是？", "恩蔭")
    ,@("栗寬(c_personid=625265的籍貫是否為源州？", "否")
    ,@("袁世承(c_personid=655454)的入仕方式是？", "恩蔭")
    ,@("楊玉山(c_personid=118508)的籍貫是？", "松江府")
    ,@("沈文明(c_personid=223855的籍貫是否為蘇州府？", "是")
    ,@("童榮(c_personid=248616)的籍貫是？", "蘭溪")
    ,@("胡永華(c_personid=648284)的籍貫是？", "懷寧")
    ,@("王稷(c_personid=189639)的籍貫是？", "河中府")
    ,@("張經(c_personid=320803的籍貫是否為揚州府？", "是")
    ,@("袁一鰲(c_personid=550662)的入仕方式是？", "舉人")
    ,@("徐紹心(c_personid=319410)的籍貫是？", "餘姚")
    ,@("夏士錦(c_personid=344320)的入仕方式是？", "進士")
    ,@("樊廷茂(c_personid=44998)的入仕方式是？", "進士")
    ,@("嚴一介(c_personid=597117)的入仕方式是？", "舉人")
    ,@("胡順(c_personid=34481)的入仕方式是？", "恩蔭")
    ,@("張經(c_personid=320803的籍貫是否為郵州？", "否")
    ,@("徐榕(c_personid=13035)的入仕方式是？", "恩蔭")
    ,@("盧庭琮(c_personid=348921的籍貫是否為武昌府？", "是")
    ,@("王綺珍(c_personid=342339的籍貫是否為臨安府？", "是")
    ,@("周熙有(c_personid=595264)的入仕方式是？", "舉人")
    ,@("趙希府(c_personid=536746)的入仕方式是？", "進士")
    ,@("戴洪緒(c_personid=613079)的入仕方式是？", "舉人")
    ,@("彭萬涵(c_personid=610647)的入仕方式是？", "舉人")
    ,@("陳鳳廷(c_personid=672711)的入仕方式是？", "舉人")
    ,@("王珝(c_personid=342024的籍貫是否為平衛？", "否")
    ,@("劉錡(c_personid=24320)的入仕方式是？", "恩蔭")
    ,@("胡汝瀾(c_personid=648294的籍貫是否為陰廳？", "否")
    ,@("譚顯相(c_personid=340256)的入仕方式是？", "進士")
    ,@("耿承祖(c_personid=105844)的入仕方式是？", "恩蔭")
    ,@("龔聞道(c_personid=206928的籍貫是否為倉州？", "否")
    ,@("齊棟(c_personid=682450)的入仕方式是？", "舉人")
    ,@("衛璜書(c_personid=655388)的籍貫是？", "趙城")
    ,@("展毓(c_personid=198547的籍貫是否為右軍都督府？", "是")
    ,@("張翕(c_personid=346964)的籍貫是？", "山東省")
    ,@("張燾(c_personid=64240)的入仕方式是？", "進士")
    ,@("胡汝瀾(c_personid=648294的籍貫是否為興安府？", "是")
    ,@("王珝(c_personid=342024的籍貫是否為後軍都督府？", "是")
    ,@("沈文明(c_personid=223855的籍貫是否為倉州？", "否")
    ,@("龔聞道(c_personid=206928的籍貫是否為蘇州府？", "是")
    ,@("蔡名載(c_personid=363279)的入仕方式是？", "進士")
    ,@("沈修文(c_personid=631676)的入仕方式是？", "舉人")
    ,@("葉得榕(c_personid=651043)的籍貫是？", "湖南省")
    ,@("唐拱(c_personid=3829)的入仕方式是？", "恩蔭")
    ,@("王廣廉(c_personid=1843)的籍貫是？", "洹水")
    ,@("栗寬(c_personid=625265的籍貫是否為大同府？", "是")
    ,@("南宮鼎(c_personid=357302)的入仕方式是？", "進士")
    ,@("王綺珍(c_personid=342339的籍貫是否為屏州？", "否")
    ,@("裴冕(c_personid=32838)的入仕方式是？", "恩蔭")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
